# Fruta / hortaliza, semanal
# Re-shuffle the weekly price records (columns D, J, K, L, M, P) across
# rows 2-22 on the active sheet, matching the reordered source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2;  D=44694; J=400; K=16000; L=17000; M=16500; P=1650},
    @{Row=3;  D=44547; J=300; K=19000; L=20000; M=19500; P=1950},
    @{Row=4;  D=44460; J=300; K=15000; L=16000; M=15500; P=1550},
    @{Row=5;  D=44428; J=300; K=15000; L=16000; M=15500; P=1550},
    @{Row=6;  D=44358; J=300; K=14000; L=15000; M=14500; P=1450},
    @{Row=7;  D=44406; J=400; K=14000; L=15000; M=14500; P=1450},
    @{Row=8;  D=44644; J=300; K=20000; L=21000; M=20500; P=2050},
    @{Row=9;  D=44441; J=300; K=15000; L=16000; M=15500; P=1550},
    @{Row=10; D=44727; J=400; K=18000; L=19000; M=18500; P=1850},
    @{Row=11; D=44204; J=400; K=10000; L=11000; M=10500; P=1050},
    @{Row=12; D=44679; J=200; K=19000; L=20000; M=19500; P=1950},
    @{Row=13; D=44160; J=360; K=10000; L=11000; M=10500; P=1050},
    @{Row=14; D=44377; J=650; K=14000; L=15000; M=14538; P=1454},
    @{Row=15; D=44263; J=300; K=15000; L=16000; M=15500; P=1550},
    @{Row=16; D=44330; J=300; K=13000; L=14000; M=13500; P=1350},
    @{Row=17; D=44218; J=320; K=10000; L=11000; M=10500; P=1050},
    @{Row=18; D=44714; J=400; K=19000; L=20000; M=19500; P=1950},
    @{Row=19; D=44580; J=200; K=18000; L=20000; M=19000; P=1900},
    @{Row=20; D=44524; J=200; K=20000; L=21000; M=20500; P=2050},
    @{Row=21; D=44291; J=200; K=13000; L=14000; M=13500; P=1350},
    @{Row=22; D=44265; J=200; K=15000; L=16000; M=15500; P=1550}
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($row, 10).Value = $r.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $r.K   # K: Precio minimo
    $ws.Cells.Item($row, 12).Value = $r.L   # L: Precio maximo
    $ws.Cells.Item($row, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $r.P   # P: Precio $/Kg
}
